$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Last Name) filled first, top to bottom
$ws.Range("B1").Value = "Last Name"
$ws.Range("B2").Value = "Baby"
$ws.Range("B3").Value = "Thomas"
$ws.Range("B4").Value = "Kumar"

# Header cells for E, D, F (in that order)
$ws.Range("E1").Value = "Confirm"
$ws.Range("D1").Value = "Password"
$ws.Range("F1").Value = "Sales percent"

# Column D values (row 2-4); E column picks up same text via value assignment
$ws.Range("D2").Value = "anu123"
$ws.Range("E2").Value = "anu123"
$ws.Range("D3").Value = "mini123"
$ws.Range("E3").Value = "mini123"
$ws.Range("D4").Value = "athi123"
$ws.Range("E4").Value = "athi123"

# Column F numeric values
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 9

# Column C (Roles) filled last, top to bottom
$ws.Range("C1").Value = "Roles"
$ws.Range("C2").Value = "Technician"
$ws.Range("C3").Value = "Executive"
$ws.Range("C4").Value = "Specialist"

$ws.Range("C6").Select()

# Touch page setup so the worksheet gets an explicit <pageSetup> element
# (portrait orientation), matching the printed-layout change in the diff.
$ws.PageSetup.Orientation = 1
